# Insert a new weekly data row before the existing row 143, shifting all
# subsequent rows down by one (dimension grows from A1:R221 to A1:R222).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(143).Insert()

$ws.Range("A143").Value2 = 9
$ws.Range("B143").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C143").Value2 = "Metropolitana"
$ws.Range("D143").Value2 = 44582
$ws.Range("E143").Value2 = 13
$ws.Range("F143").Value2 = 100112043
$ws.Range("G143").Value2 = "Pepino ensalada"
$ws.Range("H143").Value2 = "Sin especificar"
$ws.Range("I143").Value2 = "Primera"
$ws.Range("J143").Value2 = 97
$ws.Range("K143").Value2 = 12000
$ws.Range("L143").Value2 = 13000
$ws.Range("M143").Value2 = 12495
$ws.Range("N143").Value2 = "$/caja 70 unidades"
$ws.Range("O143").Value2 = "Región del Maule"
$ws.Range("P143").Value2 = 178
$ws.Range("Q143").Value2 = 70
$ws.Range("R143").Value2 = "Hortaliza"

# Match the date number formatting used by the rest of column D.
$ws.Range("D143").NumberFormat = $ws.Range("D144").NumberFormat
